$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 33
$ws_ALC.Range("H33").Value = 308.46667
$ws_ALC.Range("J33").Value = 250
$ws_ALC.Range("L33").Value = 250
$ws_ALC.Range("N33").Value = -708

# ALC row 132
$ws_ALC.Range("H132").Value = 1043.3405
$ws_ALC.Range("I132").Value = 971.04443
$ws_ALC.Range("J132").Value = 2670
$ws_ALC.Range("K132").Value = 2913.13329
$ws_ALC.Range("L132").Value = 8010
$ws_ALC.Range("M132").Value = -383.1332900000002
$ws_ALC.Range("N132").Value = -13070

# ARM row 32
$ws_ARM.Range("H32").Value = 3853.92
$ws_ARM.Range("I32").Value = 2827.5652
$ws_ARM.Range("J32").Value = 15657
$ws_ARM.Range("K32").Value = 2827.5652
$ws_ARM.Range("L32").Value = 15657
$ws_ARM.Range("M32").Value = -2540.5652
$ws_ARM.Range("N32").Value = -16231

# ARM row 74
$ws_ARM.Range("H74").Value = 7937790.5
$ws_ARM.Range("I74").Value = 9010261
$ws_ARM.Range("J74").Value = 1505.6
$ws_ARM.Range("K74").Value = 9010261
$ws_ARM.Range("L74").Value = 1505.6
$ws_ARM.Range("M74").Value = -9009387
$ws_ARM.Range("N74").Value = -3253.6

# ARM row 77
$ws_ARM.Range("H77").Value = 7937790.5
$ws_ARM.Range("I77").Value = 9010261
$ws_ARM.Range("J77").Value = 1505.6
$ws_ARM.Range("K77").Value = 45051305
$ws_ARM.Range("L77").Value = 7528
$ws_ARM.Range("M77").Value = -45046937
$ws_ARM.Range("N77").Value = -16264

# ARM row 127
$ws_ARM.Range("H127").Value = 58000
$ws_ARM.Range("J127").Value = 58000
$ws_ARM.Range("L127").Value = 58000
$ws_ARM.Range("N127").Value = -67920

# ARM row 130
$ws_ARM.Range("H130").Value = 60000
$ws_ARM.Range("J130").Value = 60000
$ws_ARM.Range("L130").Value = 60000
$ws_ARM.Range("N130").Value = -70040

# ARM row 131
$ws_ARM.Range("H131").Value = 89286
$ws_ARM.Range("J131").Value = 89286
$ws_ARM.Range("L131").Value = 89286
$ws_ARM.Range("N131").Value = -99366

# ARM row 132
$ws_ARM.Range("H132").Value = 2335.9167
$ws_ARM.Range("I132").Value = 1897.5555
$ws_ARM.Range("J132").Value = 6281.1665
$ws_ARM.Range("K132").Value = 5692.666499999999
$ws_ARM.Range("L132").Value = 18843.4995
$ws_ARM.Range("M132").Value = -3162.666499999999
$ws_ARM.Range("N132").Value = -23903.4995

# BSM row 25
$ws_BSM.Range("H25").Value = 1014
$ws_BSM.Range("I25").Value = 1014
$ws_BSM.Range("J25").Value = 0
$ws_BSM.Range("K25").Value = 1014
$ws_BSM.Range("L25").Value = 0
$ws_BSM.Range("M25").ClearContents()
$ws_BSM.Range("N25").Value = -779

# BSM row 86
$ws_BSM.Range("H86").Value = 3321.4
$ws_BSM.Range("I86").Value = 2662.4
$ws_BSM.Range("J86").Value = 5957.4
$ws_BSM.Range("K86").Value = 2662.4
$ws_BSM.Range("L86").Value = 5957.4
$ws_BSM.Range("M86").Value = -1539.4
$ws_BSM.Range("N86").Value = -8203.4

# BSM row 89
$ws_BSM.Range("H89").Value = 3321.4
$ws_BSM.Range("I89").Value = 2662.4
$ws_BSM.Range("J89").Value = 5957.4
$ws_BSM.Range("K89").Value = 13312
$ws_BSM.Range("L89").Value = 29787
$ws_BSM.Range("M89").Value = -7696
$ws_BSM.Range("N89").Value = -41019

# CRP row 2
$ws_CRP.Range("H2").Value = 2250
$ws_CRP.Range("I2").Value = 1333.3334
$ws_CRP.Range("J2").Value = 5000
$ws_CRP.Range("K2").Value = 1333.3334
$ws_CRP.Range("L2").Value = 5000
$ws_CRP.Range("M2").Value = -1220.3334
$ws_CRP.Range("N2").Value = -5226

# CRP row 22
$ws_CRP.Range("H22").Value = 1561.5625
$ws_CRP.Range("I22").Value = 324.875
$ws_CRP.Range("J22").Value = 2798.25
$ws_CRP.Range("K22").Value = 324.875
$ws_CRP.Range("L22").Value = 2798.25
$ws_CRP.Range("M22").Value = 25.125
$ws_CRP.Range("N22").Value = -3498.25

# CRP row 31
$ws_CRP.Range("H31").Value = 28626.537
$ws_CRP.Range("I31").Value = 2948.4517
$ws_CRP.Range("J31").Value = 108228.6
$ws_CRP.Range("K31").Value = 2948.4517
$ws_CRP.Range("L31").Value = 108228.6
$ws_CRP.Range("M31").Value = -2653.4517
$ws_CRP.Range("N31").Value = -108818.6

# CRP row 34
$ws_CRP.Range("H34").Value = 28626.537
$ws_CRP.Range("I34").Value = 2948.4517
$ws_CRP.Range("J34").Value = 108228.6
$ws_CRP.Range("K34").Value = 2948.4517
$ws_CRP.Range("L34").Value = 108228.6
$ws_CRP.Range("M34").Value = -2746.4517
$ws_CRP.Range("N34").Value = -108632.6

# CRP row 98
$ws_CRP.Range("H98").Value = 50000.668
$ws_CRP.Range("I98").Value = 30000
$ws_CRP.Range("J98").Value = 60001
$ws_CRP.Range("K98").Value = 30000
$ws_CRP.Range("L98").Value = 60001
$ws_CRP.Range("M98").Value = -27754
$ws_CRP.Range("N98").Value = -64493

# CUL row 4
$ws_CUL.Range("H4").Value = 11553940
$ws_CUL.Range("I4").Value = 5000095
$ws_CUL.Range("J4").Value = 23205218
$ws_CUL.Range("K4").Value = 15000285
$ws_CUL.Range("L4").Value = 69615654
$ws_CUL.Range("M4").Value = -15000173
$ws_CUL.Range("N4").Value = -69615878

# CUL row 5
$ws_CUL.Range("H5").Value = 6372.0415
$ws_CUL.Range("J5").Value = 10124.571
$ws_CUL.Range("L5").Value = 30373.713
$ws_CUL.Range("N5").Value = -30597.713

# CUL row 132
$ws_CUL.Range("H132").Value = 3958.389
$ws_CUL.Range("I132").Value = 4188.8
$ws_CUL.Range("K132").Value = 37699.2
$ws_CUL.Range("M132").Value = -35169.2

# CUL row 135
$ws_CUL.Range("H135").Value = 6372.0415
$ws_CUL.Range("J135").Value = 10124.571
$ws_CUL.Range("L135").Value = 91121.139
$ws_CUL.Range("N135").Value = -96191.139

# GSM row 18
$ws_GSM.Range("H18").Value = 21999
$ws_GSM.Range("I18").Value = 0
$ws_GSM.Range("K18").Value = 0
$ws_GSM.Range("M18").ClearContents()

# GSM row 46
$ws_GSM.Range("H46").Value = 0
$ws_GSM.Range("I46").Value = 0
$ws_GSM.Range("K46").Value = 0
$ws_GSM.Range("M46").ClearContents()

# GSM row 126
$ws_GSM.Range("H126").Value = 9671
$ws_GSM.Range("I126").Value = 4999.5
$ws_GSM.Range("K126").Value = 14998.5
$ws_GSM.Range("M126").Value = -12528.5

# GSM row 129
$ws_GSM.Range("H129").Value = 56901.332
$ws_GSM.Range("J129").Value = 59997.5
$ws_GSM.Range("L129").Value = 59997.5
$ws_GSM.Range("N129").Value = -69997.5

# LTW row 55
$ws_LTW.Range("H55").Value = 2001394
$ws_LTW.Range("J55").Value = 2574.25
$ws_LTW.Range("L55").Value = 2574.25
$ws_LTW.Range("N55").Value = -2920.25

# LTW row 124
$ws_LTW.Range("H124").Value = 30264.5
$ws_LTW.Range("J124").Value = 30264.5
$ws_LTW.Range("L124").Value = 30264.5
$ws_LTW.Range("N124").Value = -40084.5

# LTW row 125
$ws_LTW.Range("H125").Value = 50000
$ws_LTW.Range("J125").Value = 50000
$ws_LTW.Range("L125").Value = 50000
$ws_LTW.Range("N125").Value = -59840

# LTW row 129
$ws_LTW.Range("H129").Value = 50000
$ws_LTW.Range("J129").Value = 50000
$ws_LTW.Range("L129").Value = 50000
$ws_LTW.Range("N129").Value = -60000

# WVR row 107
$ws_WVR.Range("H107").Value = 1703.8823
$ws_WVR.Range("I107").Value = 1784.75
$ws_WVR.Range("K107").Value = 5354.25
$ws_WVR.Range("M107").Value = -3434.25

# WVR row 122
$ws_WVR.Range("H122").Value = 6905.795
$ws_WVR.Range("I122").Value = 2606.9167
$ws_WVR.Range("J122").Value = 13784
$ws_WVR.Range("K122").Value = 7820.750100000001
$ws_WVR.Range("L122").Value = 41352
$ws_WVR.Range("M122").Value = -5370.750100000001
$ws_WVR.Range("N122").Value = -46252

# WVR row 136
$ws_WVR.Range("H136").Value = 5006.08
$ws_WVR.Range("I136").Value = 2325.1667
$ws_WVR.Range("K136").Value = 6975.500100000001
$ws_WVR.Range("M136").Value = -4425.500100000001
